# Timing issue fix - keywords, updated tc1,2 in ubc01
#
# The CasesTab query (cell B2) is rewritten to drop the trailing
# "Cohort" column (and its OPTIONAL MATCH / coalesce expression is gone
# from the RETURN clause). The SamplesTab (B3) and FilesTab (B4) query
# text is unchanged, but the shared-string table is rebuilt/compacted
# by the engine once the old CasesTab string is no longer referenced,
# so those two cells end up pointing at different shared-string indices
# automatically - no edit needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Belgian Malinois', 'Labrador Retriever','West Highland White Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in [ 'Bladder, Prostate', 'Bladder, Urethra', 'Bladder, Urethra, Prostate', 'Urethra, Prostate'] and diag.best_response in ['Not Determined', 'Stable Disease']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# Update the CasesTab query text (B2) - removes the Cohort column.
$ws.Range("B2").Value = $casesQuery

# The wrapped-text row shrinks now that the Cohort line is gone.
$ws.Rows(2).RowHeight = 319

# Put the selection/active cell back on B2 (and scroll there), matching
# where the author was working after the edit.
$ws.Range("B2").Select()
